# Actualización automática del tracker
# Rellena resultado (G) y profit (H) para los eventos ya resueltos y añade
# la fila del nuevo evento pendiente.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper: escribe texto sin que Excel lo reinterprete (p.ej. fechas) ---
function Set-PlainText($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value2 = $text
    $range.Style = "Normal"
}

# --- 1) Resultados de partidos ya finalizados (resultado, profit) ---
$resultados = @{
    99  = @("Fallo",   -1)
    109 = @("Fallo",   -1)
    123 = @("Fallo",   -1)
    124 = @("Fallo",   -1)
    125 = @("Acierto", 1.75)
    126 = @("Fallo",   -1)
    127 = @("Acierto", 1.1)
    128 = @("Fallo",   -1)
    130 = @("Acierto", 1.5)
    132 = @("Fallo",   -1)
    136 = @("Fallo",   -1)
    137 = @("Acierto", 1.25)
    138 = @("Fallo",   -1)
    148 = @("Fallo",   -1)
}

foreach ($row in $resultados.Keys) {
    $valores = $resultados[$row]
    $resultado = $valores[0]
    $profit = $valores[1]

    $ws.Range("G$row").Value2 = $resultado
    $ws.Range("H$row").Value2 = $profit
}

# --- 2) Nueva fila con el próximo evento pendiente ---
$newRow = 151
$ws.Range("A$newRow").Value2 = 14559640
Set-PlainText $ws.Range("B$newRow") "2025-09-04"
$ws.Range("C$newRow").Value2 = "Alana Smith"
$ws.Range("D$newRow").Value2 = "Maria Kozyreva"
$ws.Range("E$newRow").Value2 = "Gana Maria Kozyreva"
$ws.Range("F$newRow").Value2 = 1.83
# "resultado" y "profit" quedan vacíos: el partido todavía no se ha jugado.
